$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.239.74'
$ws.Range("E2").Value = '  -0.87%  '

$ws.Range("D3").Value = '2.364.31'
$ws.Range("E3").Value = '  -1.52%  '

$ws.Range("E4").Value = '  +0.05%  '

$c = $ws.Range("D5")
$c.Value = '''329.57'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +3.73%  '

$c = $ws.Range("D6")
$c.Value = '''107.30'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -6.71%  '

$c = $ws.Range("D7")
$c.Value = '''0.636'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -0.30%  '

$ws.Range("E8").Value = '  +0.05%  '

$c = $ws.Range("D9")
$c.Value = '''0.613'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -2.59%  '

$c = $ws.Range("D10")
$c.Value = '''41.02'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -3.99%  '

$ws.Range("E11").Value = '  -1.60%  '

$c = $ws.Range("D12")
$c.Value = '''8.45'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -3.61%  '

$ws.Range("E13").Value = '  -0.39%  '

$c = $ws.Range("D14")
$c.Value = '''0.977'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -3.60%  '

$ws.Range("D15").Value = '2.724.02'
$ws.Range("E15").Value = '  -1.64%  '

$c = $ws.Range("D16")
$c.Value = '''15.37'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -3.93%  '

$ws.Range("D17").Value = '2.355.54'
$ws.Range("E17").Value = '  -2.05%  '

$ws.Range("D18").Value = '45.207.91'
$ws.Range("E18").Value = '  -1.01%  '

$c = $ws.Range("D19")
$c.Value = '''15.30'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +12.43%  '

$c = $ws.Range("D20")
$c.Value = '''7.29'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -3.28%  '

$ws.Range("E21").Value = '  -2.21%  '

$c = $ws.Range("D22")
$c.Value = '''3.65'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +2.49%  '

$c = $ws.Range("D23")
$c.Value = '''73.04'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -2.70%  '

$c = $ws.Range("D24")
$c.Value = '''260.12'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -2.02%  '

$c = $ws.Range("D25")
$c.Value = '''2.30'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -3.48%  '

$ws.Range("E26").Value = '  -0.18%  '

$c = $ws.Range("D27")
$c.Value = '''11.33'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -0.66%  '

$ws.Range("E28").Value = '  -2.70%  '

$c = $ws.Range("D29")
$c.Value = '''2.29'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -2.38%  '

$c = $ws.Range("D30")
$c.Value = '''22.32'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -2.12%  '

$ws.Range("E31").Value = '  -3.48%  '

$c = $ws.Range("D32")
$c.Value = '''36.89'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -9.26%  '

$c = $ws.Range("D33")
$c.Value = '''167.28'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -3.40%  '

$c = $ws.Range("D34")
$c.Value = '''2.83'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -3.41%  '

$c = $ws.Range("D36")
$c.Value = '''3.27'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +4.44%  '

$ws.Range("E37").Value = '  -2.47%  '

$c = $ws.Range("D38")
$c.Value = '''4.73'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -5.81%  '

$c = $ws.Range("D39")
$c.Value = '''1.95'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +9.71%  '

$c = $ws.Range("D40")
$c.Value = '''4.00'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -6.01%  '

$ws.Range("E41").Value = '  -3.36%  '

$c = $ws.Range("D42")
$c.Value = '''97.22'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -3.03%  '

$c = $ws.Range("D43")
$c.Value = '''70.10'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -2.55%  '

$ws.Range("D44").Value = '1.882.77'
$ws.Range("E44").Value = '  +13.81%  '

$c = $ws.Range("D45")
$c.Value = '''0.229'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -4.94%  '

$c = $ws.Range("D46")
$c.Value = '''6.08'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +4.28%  '

$c = $ws.Range("D47")
$c.Value = '''12.91'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -7.13%  '

$ws.Range("B48").Value = 'FirstDigitalUSD'
$ws.Range("C48").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$c = $ws.Range("D48")
$c.Value = '''1.00'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +0.33%  '

$ws.Range("B49").Value = 'ordi'
$ws.Range("C49").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$c = $ws.Range("D49")
$c.Value = '''85.55'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -1.43%  '

$c = $ws.Range("D50")
$c.Value = '''112.36'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -3.65%  '

$c = $ws.Range("D51")
$c.Value = '''9.29'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -2.84%  '
